# Update row 4 (青田·未闻展名国漫嘉年华) ticket stats on both the
# "展览" sheet and the "全部类型" sheet: 想去人数 (F4) 4 -> 8,
# 最低票价 (G4) 39.9 -> 45.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 8
    $ws.Range("G4").Value = 45
}
